# Merge tables: insert a new "img" column (column D) into the sheet,
# shifting the existing reproduction/appearance/key_features/img columns
# one position to the right (D->E, E->F, F->G, G->H), and repair the
# hyperlinks that live in the (now shifted) last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D (the existing "habitat" column
#    stays at C, everything from "reproduction" onward shifts right by one).
$ws.Columns.Item(4).Insert()

# The inserted column picks up the width of its left neighbour (column C),
# matching Excel's normal "insert column" behaviour.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# 2. The new column D reuses the "img" header text that already lives in
#    the (now shifted) last column H2.
$ws.Range("D2").Value = $ws.Range("H2").Value2

# 3. The hyperlinks that used to anchor on column G now need to point at
#    column H. This runtime does not update hyperlink anchors when columns
#    are inserted, so rebuild them explicitly in their original order.
$hyperlinkData = @(
  @{Row=3;  Url='https://th.bing.com/th/id/R.fa385930af77ee7a1225f7093b09569b?rik=4W4tiXRbyQ7lIg&riu=http%3a%2f%2fmedia-cache-ak0.pinimg.com%2f1200x%2f2c%2f9c%2ffd%2f2c9cfd90b6bdc31da72013bd004a3402.jpg&ehk=ZdvCDcvV%2fzXJWIszU8VKIw%2bZeHylQfhR80tGsCThL2Q%3d&risl=&pid=ImgRaw&r=0'; Display='https://th.bing.com/th/id/R.fa385930af77ee7a1225f7093b09569b?rik=4W4tiXRbyQ7lIg&riu=http%3a%2f%2fmedia-cache-ak0.pinimg.com%2f1200x%2f2c%2f9c%2ffd%2f2c9cfd90b6bdc31da72013bd004a3402.jpg&ehk=ZdvCDcvV%2fzXJWIszU8VKIw%2bZeHylQfhR80tGsCThL2Q%3d&risl=&pid=ImgRaw&r=0'}
  @{Row=4;  Url='https://wallpapercave.com/w/wp7419706'; Display=''}
  @{Row=5;  Url='https://wallpapercave.com/w/wp12011971'; Display=''}
  @{Row=6;  Url='https://gameluster.com/wp-content/uploads/2023/02/Unicorn-Patronus.jpg'; Display=''}
  @{Row=7;  Url='https://i.pinimg.com/originals/0a/2d/9d/0a2d9d0e1d4d556ce76f4a254ed410f6.jpg'; Display=''}
  @{Row=8;  Url='https://pm1.narvii.com/6774/bf0a848bdd54c9ff23a111310f6d05628c23bf3dv2_hq.jpg'; Display=''}
  @{Row=9;  Url='https://i.insider.com/582d0dc165edfe1b008b46c7?width=1300&format=jpeg&auto=webp'; Display=''}
  @{Row=10; Url='https://tse1.mm.bing.net/th?id=OIP.pO1q0uL0as8jujGz6v45cgHaEM&pid=ImgDet&rs=1'; Display=''}
  @{Row=11; Url='https://th.bing.com/th/id/OIP.OZz36ap2VZ9FV8hQ_qKaAAHaMN?pid=ImgDet&rs=1'; Display=''}
  @{Row=13; Url='https://thehumblefabulist.files.wordpress.com/2018/10/20170714_121019.jpg?w=1024&h=806'; Display=''}
  @{Row=12; Url='https://e1.pxfuel.com/desktop-wallpaper/742/998/desktop-wallpaper-harry-potter-book-harry-potter-dragon.jpg'; Display=''}
)

# Remove every existing (stale, G-column) hyperlink in one shot.
$ws.Hyperlinks.Delete()

foreach ($entry in $hyperlinkData) {
  $target = $ws.Cells.Item($entry.Row, 8)
  if ($entry.Display -ne '') {
    $ws.Hyperlinks.Add($target, $entry.Url, "", "", $entry.Display)
  } else {
    $ws.Hyperlinks.Add($target, $entry.Url)
  }
}

# 4. Match the author's final cursor position.
[void]$ws.Range("C19").Select()
